# Applies the "cleaned data" edit: the Authors values in E2 and E3 get one
# additional padding space inserted between each comma-separated author
# entry (matching the new shared-string entries added at the end of
# sharedStrings.xml, replacing the references to the old entries).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the accented character (é, U+00E9) safely regardless of file encoding.
$eAcute = [char]0x00E9

$pad18 = "                  "

$e2 = '[Anthony V%Das%NULL%1,' + $pad18 + 'Padmaja K%Rani%NULL%1,' + $pad18 + 'Pravin K%Vaddavalli%NULL%1]'

$e3 = '[Gagan%Kalra%NULL%1,' + $pad18 + 'Andrew M.%Williams%NULL%1,' + $pad18 + 'Patrick W.%Commiskey%NULL%1,' + $pad18 + 'Eve M. R.%Bowers%NULL%1,' + $pad18 + 'Tadhg%Schempf%NULL%1,' + $pad18 + 'Jos' + $eAcute + '-Alain%Sahel%NULL%1,' + $pad18 + 'Evan L.%Waxman%waxmane@upmc.edu%1,' + $pad18 + 'Roxana%Fu%fur3@upmc.edu%1]'

$ws.Range("E2").Value = $e2
$ws.Range("E3").Value = $e3
